$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Save the "closing row" (bottom-border) formatting that currently lives
#    on row 18 into a scratch row far below the used range, so we can re-
#    apply it later to the new last data row (44) after the insert shifts
#    everything down.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B100:J100").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Make room for the additional worker/period rows. Row 18 (first data row
#    with the "closing" style) stays put; we insert 26 new rows right after
#    it (rows 19-44) which pushes the old trailing rows (23 "___" and 24
#    "NOMBRE...") down to rows 49 and 50.
# ---------------------------------------------------------------------------
$ws.Rows("19:44").Insert()

# ---------------------------------------------------------------------------
# 3) Row 18 is no longer the last data row, so give it the regular
#    (no-bottom-border) formatting, copying it from row 17 which still has
#    the correct "middle" style.
# ---------------------------------------------------------------------------
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Row 44 is now the new last data row -> give it the "closing" style we
#    stashed away in step 1 (now located at row 126 because of the insert).
# ---------------------------------------------------------------------------
$ws.Range("B126:J126").Copy()
$ws.Range("B44:J44").PasteSpecial(-4122)

# Remove the scratch row used to stash the formatting.
$ws.Rows("126:126").Delete()

# ---------------------------------------------------------------------------
# 4b) Every freshly inserted row (19-43) currently has generic/default
#     formatting; give them the same "middle" row style now sitting on
#     row 18.
# ---------------------------------------------------------------------------
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J43").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 5) Populate the detail rows (16-44) with the updated worker / mora data.
# ---------------------------------------------------------------------------
$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "73213860"
$ws.Cells.Item(16, 4).Value = "ARIEL JAMID TOSCANO MENDOZA"
$ws.Cells.Item(16, 5).Value = "2201"
$ws.Cells.Item(16, 6).Value = 36341
$ws.Cells.Item(16, 7).Value = 908526

$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "73213860"
$ws.Cells.Item(17, 4).Value = "ARIEL JAMID TOSCANO MENDOZA"
$ws.Cells.Item(17, 5).Value = "2112"
$ws.Cells.Item(17, 6).Value = 36341
$ws.Cells.Item(17, 7).Value = 908526

$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "1143332016"
$ws.Cells.Item(18, 4).Value = "JOSE OSCAR ACOSTA HOYOS"
$ws.Cells.Item(18, 5).Value = "2106"
$ws.Cells.Item(18, 6).Value = 1211
$ws.Cells.Item(18, 7).Value = 908526

$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "12598379"
$ws.Cells.Item(19, 4).Value = "ILMER IVAN PASSO PUELLO"
$ws.Cells.Item(19, 5).Value = "2507"
$ws.Cells.Item(19, 6).Value = 56940
$ws.Cells.Item(19, 7).Value = 908526

$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "12598379"
$ws.Cells.Item(20, 4).Value = "ILMER IVAN PASSO PUELLO"
$ws.Cells.Item(20, 5).Value = "2506"
$ws.Cells.Item(20, 6).Value = 56940
$ws.Cells.Item(20, 7).Value = 908526

$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "9148712"
$ws.Cells.Item(21, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(21, 5).Value = "2507"
$ws.Cells.Item(21, 6).Value = 36341
$ws.Cells.Item(21, 7).Value = 908526

$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "9148712"
$ws.Cells.Item(22, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(22, 5).Value = "2506"
$ws.Cells.Item(22, 6).Value = 36341
$ws.Cells.Item(22, 7).Value = 908526

$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "9148712"
$ws.Cells.Item(23, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(23, 5).Value = "2505"
$ws.Cells.Item(23, 6).Value = 36341
$ws.Cells.Item(23, 7).Value = 908526

$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "9148712"
$ws.Cells.Item(24, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(24, 5).Value = "2504"
$ws.Cells.Item(24, 6).Value = 36341
$ws.Cells.Item(24, 7).Value = 908526

$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "9148712"
$ws.Cells.Item(25, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(25, 5).Value = "2503"
$ws.Cells.Item(25, 6).Value = 36341
$ws.Cells.Item(25, 7).Value = 908526

$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "9148712"
$ws.Cells.Item(26, 4).Value = "LUIS CARLOS BLANCO VERGARA"
$ws.Cells.Item(26, 5).Value = "2502"
$ws.Cells.Item(26, 6).Value = 36341
$ws.Cells.Item(26, 7).Value = 908526

$ws.Cells.Item(27, 2).Value = "CC"
$ws.Cells.Item(27, 3).Value = "12589228"
$ws.Cells.Item(27, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(27, 5).Value = "2507"
$ws.Cells.Item(27, 6).Value = 36341
$ws.Cells.Item(27, 7).Value = 877803

$ws.Cells.Item(28, 2).Value = "CC"
$ws.Cells.Item(28, 3).Value = "12589228"
$ws.Cells.Item(28, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(28, 5).Value = "2506"
$ws.Cells.Item(28, 6).Value = 36341
$ws.Cells.Item(28, 7).Value = 877803

$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "12589228"
$ws.Cells.Item(29, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(29, 5).Value = "2505"
$ws.Cells.Item(29, 6).Value = 36341
$ws.Cells.Item(29, 7).Value = 877803

$ws.Cells.Item(30, 2).Value = "CC"
$ws.Cells.Item(30, 3).Value = "12589228"
$ws.Cells.Item(30, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(30, 5).Value = "2504"
$ws.Cells.Item(30, 6).Value = 36341
$ws.Cells.Item(30, 7).Value = 877803

$ws.Cells.Item(31, 2).Value = "CC"
$ws.Cells.Item(31, 3).Value = "12589228"
$ws.Cells.Item(31, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(31, 5).Value = "2503"
$ws.Cells.Item(31, 6).Value = 36341
$ws.Cells.Item(31, 7).Value = 877803

$ws.Cells.Item(32, 2).Value = "CC"
$ws.Cells.Item(32, 3).Value = "12589228"
$ws.Cells.Item(32, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(32, 5).Value = "2502"
$ws.Cells.Item(32, 6).Value = 36341
$ws.Cells.Item(32, 7).Value = 877803

$ws.Cells.Item(33, 2).Value = "CC"
$ws.Cells.Item(33, 3).Value = "12589228"
$ws.Cells.Item(33, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(33, 5).Value = "2109"
$ws.Cells.Item(33, 6).Value = 36341
$ws.Cells.Item(33, 7).Value = 877803

$ws.Cells.Item(34, 2).Value = "CC"
$ws.Cells.Item(34, 3).Value = "12589228"
$ws.Cells.Item(34, 4).Value = "JUAN IRABA GUZMAN DE HOYOS"
$ws.Cells.Item(34, 5).Value = "2108"
$ws.Cells.Item(34, 6).Value = 36341
$ws.Cells.Item(34, 7).Value = 877803

$ws.Cells.Item(35, 2).Value = "CC"
$ws.Cells.Item(35, 3).Value = "1047453515"
$ws.Cells.Item(35, 4).Value = "JOSE JADER VERGARA BARBOZA"
$ws.Cells.Item(35, 5).Value = "2502"
$ws.Cells.Item(35, 6).Value = 52000
$ws.Cells.Item(35, 7).Value = 1300000

$ws.Cells.Item(36, 2).Value = "CC"
$ws.Cells.Item(36, 3).Value = "1047445983"
$ws.Cells.Item(36, 4).Value = "DAYANA PAOLA ALTAMAR DIAZ"
$ws.Cells.Item(36, 5).Value = "2507"
$ws.Cells.Item(36, 6).Value = 56940
$ws.Cells.Item(36, 7).Value = 1423500

$ws.Cells.Item(37, 2).Value = "CC"
$ws.Cells.Item(37, 3).Value = "1047445983"
$ws.Cells.Item(37, 4).Value = "DAYANA PAOLA ALTAMAR DIAZ"
$ws.Cells.Item(37, 5).Value = "2506"
$ws.Cells.Item(37, 6).Value = 56940
$ws.Cells.Item(37, 7).Value = 1423500

$ws.Cells.Item(38, 2).Value = "CC"
$ws.Cells.Item(38, 3).Value = "1007254953"
$ws.Cells.Item(38, 4).Value = "IVANNA PAOLA PASSO CORREA"
$ws.Cells.Item(38, 5).Value = "2502"
$ws.Cells.Item(38, 6).Value = 52000
$ws.Cells.Item(38, 7).Value = 908526

$ws.Cells.Item(39, 2).Value = "CC"
$ws.Cells.Item(39, 3).Value = "9094562"
$ws.Cells.Item(39, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(39, 5).Value = "2507"
$ws.Cells.Item(39, 6).Value = 36341
$ws.Cells.Item(39, 7).Value = 908526

$ws.Cells.Item(40, 2).Value = "CC"
$ws.Cells.Item(40, 3).Value = "9094562"
$ws.Cells.Item(40, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(40, 5).Value = "2506"
$ws.Cells.Item(40, 6).Value = 36341
$ws.Cells.Item(40, 7).Value = 908526

$ws.Cells.Item(41, 2).Value = "CC"
$ws.Cells.Item(41, 3).Value = "9094562"
$ws.Cells.Item(41, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(41, 5).Value = "2505"
$ws.Cells.Item(41, 6).Value = 36341
$ws.Cells.Item(41, 7).Value = 908526

$ws.Cells.Item(42, 2).Value = "CC"
$ws.Cells.Item(42, 3).Value = "9094562"
$ws.Cells.Item(42, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(42, 5).Value = "2504"
$ws.Cells.Item(42, 6).Value = 36341
$ws.Cells.Item(42, 7).Value = 908526

$ws.Cells.Item(43, 2).Value = "CC"
$ws.Cells.Item(43, 3).Value = "9094562"
$ws.Cells.Item(43, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(43, 5).Value = "2503"
$ws.Cells.Item(43, 6).Value = 36341
$ws.Cells.Item(43, 7).Value = 908526

$ws.Cells.Item(44, 2).Value = "CC"
$ws.Cells.Item(44, 3).Value = "9094562"
$ws.Cells.Item(44, 4).Value = "FEDERICO ANTONIO LARA BERRIO"
$ws.Cells.Item(44, 5).Value = "2502"
$ws.Cells.Item(44, 6).Value = 36341
$ws.Cells.Item(44, 7).Value = 908526

# ---------------------------------------------------------------------------
# 6) Update the summary figures: total overdue value, worker count and
#    period count.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1132473
$ws.Range("C13").Value = 9
$ws.Range("F13").Value = 11

# ---------------------------------------------------------------------------
# 7) Column D (worker name) needs to be a bit wider to fit the longest new
#    name ("FEDERICO ANTONIO LARA BERRIO").
# ---------------------------------------------------------------------------
$ws.Columns("D:D").AutoFit()
